# Commit: "Changed format in several Excel files. Adjusted portfolios
# classification chart"
#
# The bulk of the underlying diff (fileVersion/rupBuild, x15ac:absPath,
# xr:revisionPtr, bookViews window geometry, theme display name, and the
# sheetFormatPr baseColWidth/defaultColWidth attribute swap) is just Excel
# re-saving the workbook on a different machine/build - environment
# metadata that isn't part of the document's actual content model. The one
# real, content-level edit captured by the diff is the worksheet rename:
#   "Hoja1" -> "Data"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Name = "Data"
